$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "48.134.89"
$ws.Range("E2").Value = "  +2.33%  "

# Row 3
$ws.Range("D3").Value = "2.516.38"
$ws.Range("E3").Value = "  +1.38%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.07"
$ws.Range("E5").Value = "  +0.72%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.94"
$ws.Range("E6").Value = "  +0.87%  "

# Row 7
$ws.Range("E7").Value = "  +2.18%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.553"
$ws.Range("E9").Value = "  +4.23%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.37"
$ws.Range("E10").Value = "  +4.81%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.55"
$ws.Range("E11").Value = "  +13.42%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0823"
$ws.Range("E12").Value = "  +2.10%  "

# Row 13
$ws.Range("E13").Value = "  +1.26%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.26"
$ws.Range("E14").Value = "  +2.44%  "

# Row 15
$ws.Range("D15").Value = "2.913.41"
$ws.Range("E15").Value = "  +1.25%  "

# Row 16
$ws.Range("D16").Value = "2.516.75"
$ws.Range("E16").Value = "  +1.13%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.853"
$ws.Range("E17").Value = "  +1.40%  "

# Row 18
$ws.Range("D18").Value = "47.970.15"
$ws.Range("E18").Value = "  +2.02%  "

# Row 19
$ws.Range("E19").Value = "  +4.64%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.62"
$ws.Range("E20").Value = "  +0.87%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0946"
$ws.Range("E21").Value = "  +1.65%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.71"
$ws.Range("E22").Value = "  -1.67%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.00"
$ws.Range("E23").Value = "  +1.98%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.87"
$ws.Range("E24").Value = "  +7.78%  "

# Row 25
$ws.Range("E25").Value = "  +0.81%  "

# Row 26
$ws.Range("E26").Value = "  -0.29%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.05"
$ws.Range("E27").Value = "  +1.90%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("E28").Value = "  +1.04%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.144"
$ws.Range("E29").Value = "  +1.99%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.20"
$ws.Range("E30").Value = "  -3.08%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.19"
$ws.Range("E31").Value = "  +4.01%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.75"
$ws.Range("E32").Value = "  +0.20%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.74"
$ws.Range("E33").Value = "  -0.89%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.40"
$ws.Range("E34").Value = "  +1.49%  "

# Row 35
$ws.Range("E35").Value = "  -0.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0789"
$ws.Range("E36").Value = "  +1.22%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.98"
$ws.Range("E37").Value = "  +1.69%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.71"
$ws.Range("E38").Value = "  +1.92%  "

# Row 39
$ws.Range("E39").Value = "  +1.98%  "

# Row 40
$ws.Range("E40").Value = "  +0.68%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.17"
$ws.Range("E41").Value = "  +3.77%  "

# Row 42
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.20"
$ws.Range("E42").Value = "  -0.76%  "

# Row 43
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.55"
$ws.Range("E43").Value = "  +0.30%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0300"

# Row 45
$ws.Range("D45").Value = "2.014.92"
$ws.Range("E45").Value = "  +2.02%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.16"
$ws.Range("E46").Value = "  +5.49%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.91"
$ws.Range("E47").Value = "  +9.37%  "

# Row 48
$ws.Range("E48").Value = "  +2.44%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.12"
$ws.Range("E49").Value = "  +0.85%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.20"
$ws.Range("E50").Value = "  +1.73%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.77"
$ws.Range("E51").Value = "  +3.07%  "
